$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh naive forecaster QoQ forecast-error stats for ifo GDP component analysis (Q0-Q8)
$ws.Range("B2").Value = -0.07937728594090611
$ws.Range("C2").Value = 1.316068219616539
$ws.Range("D2").Value = 9.638393705627079
$ws.Range("E2").Value = 3.104576252184359
$ws.Range("F2").Value = 3.13384069789644
$ws.Range("G2").Value = 52

$ws.Range("B3").Value = -0.02448439541659127
$ws.Range("C3").Value = 1.462623545219845
$ws.Range("D3").Value = 9.051394970963811
$ws.Range("E3").Value = 3.008553634383773
$ws.Range("F3").Value = 3.038389605539716
$ws.Range("G3").Value = 51

$ws.Range("B4").Value = -0.08143474300426477
$ws.Range("C4").Value = 1.364967268314127
$ws.Range("D4").Value = 6.831940587747885
$ws.Range("E4").Value = 2.613798115338651
$ws.Range("F4").Value = 2.639053048004792
$ws.Range("G4").Value = 50

$ws.Range("B5").Value = -0.001877683171190082
$ws.Range("C5").Value = 1.602456891704275
$ws.Range("D5").Value = 9.4526324803627
$ws.Range("E5").Value = 3.074513372936065
$ws.Range("F5").Value = 3.106373886794302
$ws.Range("G5").Value = 49

$ws.Range("B6").Value = -0.1039814330460372
$ws.Range("C6").Value = 1.471754775803563
$ws.Range("D6").Value = 8.504042060525727
$ws.Range("E6").Value = 2.91616907269207
$ws.Range("F6").Value = 2.945154824746354
$ws.Range("G6").Value = 48

$ws.Range("B7").Value = -0.1015771582092904
$ws.Range("C7").Value = 1.769024061707294
$ws.Range("D7").Value = 10.29407639076258
$ws.Range("E7").Value = 3.208438310262889
$ws.Range("F7").Value = 3.252319194407661
$ws.Range("G7").Value = 36

$ws.Range("B8").Value = -0.1097886217413366
$ws.Range("C8").Value = 1.815971081873004
$ws.Range("D8").Value = 10.93971098436304
$ws.Range("E8").Value = 3.307523391355387
$ws.Range("F8").Value = 3.353961707055734
$ws.Range("G8").Value = 35

$ws.Range("B9").Value = 0.007295727682830271
$ws.Range("C9").Value = 2.609854038779015
$ws.Range("D9").Value = 19.14435649662439
$ws.Range("E9").Value = 4.375426435974485
$ws.Range("F9").Value = 4.50227040007513
$ws.Range("G9").Value = 18

$ws.Range("B10").Value = -0.6545140871090607
$ws.Range("C10").Value = 2.980052884870046
$ws.Range("D10").Value = 23.04619067002914
$ws.Range("E10").Value = 4.800644818149864
$ws.Range("F10").Value = 4.987943682298893
$ws.Range("G10").Value = 11
